$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 421, shifting existing rows 421..502 down to 422..503.
$ws.Rows.Item(421).Insert()

# Populate the newly inserted row 421 with the new record's data.
# Columns A, B, C, E, F, G, H, I, Q, R keep the same values as the (now shifted-down)
# row 422, which is why Excel's Insert already carries matching formatting through;
# we only need to set the specific data cells for the new record.
$ws.Cells.Item(421, 1).Value = 6
$ws.Cells.Item(421, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(421, 3).Value = "Metropolitana"
$ws.Cells.Item(421, 4).Value = 44694
$ws.Cells.Item(421, 5).Value = 13
$ws.Cells.Item(421, 6).Value = 100112044
$ws.Cells.Item(421, 7).Value = "Perejil"
$ws.Cells.Item(421, 8).Value = "Sin especificar"
$ws.Cells.Item(421, 9).Value = "Primera"
$ws.Cells.Item(421, 10).Value = 220
$ws.Cells.Item(421, 11).Value = 11000
$ws.Cells.Item(421, 12).Value = 12000
$ws.Cells.Item(421, 13).Value = 11409
$ws.Cells.Item(421, 14).Value = "$/docena de atados"
$ws.Cells.Item(421, 15).Value = "Región Metropolitana"
$ws.Cells.Item(421, 16).Value = 3803
$ws.Cells.Item(421, 17).Value = 3
$ws.Cells.Item(421, 18).Value = "Hortaliza"
